# Fruta / hortaliza, semanal
# Insert two new weekly price records for "Sweet Heart" cherries at
# Terminal Hortofrutícola Agro Chillán (row 61/62), pushing the existing
# data (old rows 61-152) down to rows 63-154.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 61; everything that used
# to start at row 61 now starts at row 63.
$ws.Rows("61:62").Insert()

# New row 61: Sweet Heart / Primera
$ws.Range("A61").Value = 7
$ws.Range("B61").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C61").Value = "Ñuble"
$ws.Range("D61").Value = 44930
$ws.Range("E61").Value = 16
$ws.Range("F61").Value = "Fruta"
$ws.Range("G61").Value = 100103
$ws.Range("H61").Value = "Frutos de hueso (carozo)"
$ws.Range("I61").Value = 100103001
$ws.Range("J61").Value = "Cereza"
$ws.Range("K61").Value = "Sweet Heart"
$ws.Range("L61").Value = "Primera"
$ws.Range("M61").Value = 80
$ws.Range("N61").Value = 5000
$ws.Range("O61").Value = 5000
$ws.Range("P61").Value = 5000
$ws.Range("Q61").Value = "$/bandeja 10 kilos"
$ws.Range("R61").Value = "Quillón"
$ws.Range("S61").Value = 500
$ws.Range("T61").Value = 10

# New row 62: Sweet Heart / Segunda
$ws.Range("A62").Value = 7
$ws.Range("B62").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C62").Value = "Ñuble"
$ws.Range("D62").Value = 44930
$ws.Range("E62").Value = 16
$ws.Range("F62").Value = "Fruta"
$ws.Range("G62").Value = 100103
$ws.Range("H62").Value = "Frutos de hueso (carozo)"
$ws.Range("I62").Value = 100103001
$ws.Range("J62").Value = "Cereza"
$ws.Range("K62").Value = "Sweet Heart"
$ws.Range("L62").Value = "Segunda"
$ws.Range("M62").Value = 80
$ws.Range("N62").Value = 4000
$ws.Range("O62").Value = 4000
$ws.Range("P62").Value = 4000
$ws.Range("Q62").Value = "$/bandeja 10 kilos"
$ws.Range("R62").Value = "Quillón"
$ws.Range("S62").Value = 400
$ws.Range("T62").Value = 10
